$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = "  -0.41%  "
    "E3" = "  -0.66%  "
    "E5" = "  -0.09%  "
    "E6" = "  +0.74%  "
    "E7" = "  +0.04%  "
    "E8" = "  +0.87%  "
    "E9" = "  +4.95%  "
    "E10" = "  -1.25%  "
    "E11" = "  +2.08%  "
    "E12" = "  -0.76%  "
    "E13" = "  +0.02%  "
    "E14" = "  -5.58%  "
    "E15" = "  -0.82%  "
    "E16" = "  -0.03%  "
    "E17" = "  -0.44%  "
    "E18" = "  -0.21%  "
    "E19" = "  +1.82%  "
    "E20" = "  -3.02%  "
    "E21" = "  -0.89%  "
    "E22" = "  -0.14%  "
    "E23" = "  +0.25%  "
    "E24" = "  +0.10%  "
    "E25" = "  -0.93%  "
    "E26" = "  -0.24%  "
    "E27" = "  +1.83%  "
    "E28" = "  -0.63%  "
    "E29" = "  -0.28%  "
    "E30" = "  -2.22%  "
    "E31" = "  -1.50%  "
    "E32" = "  -0.01%  "
    "E33" = "  -5.15%  "
    "E34" = "  -2.49%  "
    "E35" = "  +1.79%  "
    "E36" = "  +3.72%  "
    "E37" = "  +1.71%  "
    "E38" = "  -2.75%  "
    "E39" = "  -0.77%  "
    "E40" = "  -0.86%  "
    "E41" = "  +0.40%  "
    "E42" = "  -1.97%  "
    "E43" = "  +0.87%  "
    "E44" = "  -2.95%  "
    "E45" = "  -1.43%  "
    "E46" = "  -1.74%  "
    "E47" = "  -0.78%  "
    "E48" = "  +1.62%  "
    "E49" = "  -0.41%  "
    "E50" = "  -4.06%  "
    "E51" = "  -0.02%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Updated $($updates.Count) cells"
